$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) for the two events that appear on
# both the "展览" sheet and the "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 183
    $ws.Range("F4").Value = 135
}
